# Monte Carlo / Sensitivity variable table updates
# Adds newly-supported sensitivity/Monte-Carlo parameters (CostInvestIncr,
# CostFixedIncr, CapacityCredit, StorageDuration, HeatRate, LastBuild,
# FirstBuild, Retirement) to the PowerPlants, Fuels and Connections
# variable tables.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# PowerPlants sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("PowerPlants")

$powerPlants = @(
    @("CapacityCredit",   "N"),
    @("StorageDuration",  "N"),
    @("CapacityFactor",   "N"),
    @("Efficiency",       "N"),
    @("HeatRate",         "N"),
    @("ExpectedLifetime", "N"),
    @("CostInvest",       "Y"),
    @("CostInvestIncr",   "N"),
    @("CostFixed",        "N"),
    @("CostFixedIncr",    "N"),
    @("CostVariable",     "N"),
    @("CostVariableIncr", "N"),
    @("DiscountRate",     "N"),
    @("RampRate",         "N"),
    @("MaxCapacity",      "N"),
    @("MaxActivity",      "N"),
    @("FirstBuild",       "N"),
    @("LastBuild",        "N")
)

$r = 2
foreach ($row in $powerPlants) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $r = $r + 1
}

# Rows 20-33 remain blank placeholder rows (extends the table down to B33)
for ($r = 20; $r -le 33; $r++) {
    if ($ws.Cells.Item($r, 1).Value -eq $null) {
        $ws.Cells.Item($r, 1).Value = ""
    }
    if ($ws.Cells.Item($r, 2).Value -eq $null) {
        $ws.Cells.Item($r, 2).Value = ""
    }
}

$ws.Range("A23").Select()

# ---------------------------------------------------------------------
# Fuels sheet
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Fuels")
$ws3.Activate()

$fuels = @(
    @("CostVariable",     "N"),
    @("CostVariableIncr", "N"),
    @("CostInvest",       "N"),
    @("CostInvestIncr",   "N"),
    @("EmissionActivity", "N"),
    @("Lifetime",         "N"),
    @("Retirement",       "N"),
    @("FirstBuild",       "N"),
    @("LastBuild",        "N"),
    @("MaxActivity",      "N")
)

$r = 2
foreach ($row in $fuels) {
    $ws3.Cells.Item($r, 1).Value = $row[0]
    $ws3.Cells.Item($r, 2).Value = $row[1]
    $r = $r + 1
}

$ws3.Range("B12").Select()

# ---------------------------------------------------------------------
# Connections sheet
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Connections")
$ws4.Activate()

$connections = @(
    @("Loss",             "N"),
    @("CostVariable",     "N"),
    @("CostVariableIncr", "N"),
    @("CostInvest",       "N"),
    @("CostInvestIncr",   "N"),
    @("EmissionActivity", "N"),
    @("Lifetime",         "N"),
    @("FirstBuild",       "N"),
    @("LastBuild",        "N")
)

$r = 2
foreach ($row in $connections) {
    $ws4.Cells.Item($r, 1).Value = $row[0]
    $ws4.Cells.Item($r, 2).Value = $row[1]
    $r = $r + 1
}

$ws4.Columns.Item(1).AutoFit()
$ws4.Range("C8").Select()

# ---------------------------------------------------------------------
# Restore the Globals sheet as the active tab (matches original workbook
# view state, which this edit does not otherwise touch)
# ---------------------------------------------------------------------
$wsGlobals = $wb.Worksheets.Item("Globals")
$wsGlobals.Activate()
$wsGlobals.Range("B3").Select()
